$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add new header values P1, Q1 (copy style from O1 which has s="1") ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: swap I/K and M/O values ---
# Before: I=1, K=2, M=1, O=2  ->  After: I=2, K=1, M=2, O=1
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# --- Rows 2-25: add new columns P and Q, both value 2 ---
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
